$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update text values in column B (note text, suffix with IBAN fragment)
$ws.Range("B2").Value = "Anticipo Ft. Estero B.Cred.T. (IT74*680)"
$ws.Range("B3").Value = "Anticipo Ft. Italia B.Pop.Soft. (IT15*456)"

# Rename the journal reference used throughout column J, and drop the
# stray direct formatting those cells carried (reverts them to the
# worksheet's default/general style).
$ws.Range("J2:J7").Clear()
$ws.Range("J2:J7").Value = "z0bug.jou_misc"

# Widen column B and bump the sheet's default column width slightly
$ws.Columns.Item(2).ColumnWidth = 41.25
$ws.StandardWidth = 11.70703125

# Move the active selection to B2
$ws.Range("B2").Select()
